# Weekly refresh of Hortaliza / Macroferia Regional de Talca - Haba data.
# Existing rows 2-31 get updated Fecha/Volumen/Precio/Origen values for the
# new week, and one additional weekly observation is appended as row 32.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing data rows (2-31) ---
$ws.Range("D2").Value = 44162
$ws.Range("J2").Value = 200
$ws.Range("O2").Value = 'Región del Maule'
$ws.Range("D3").Value = 44396
$ws.Range("K3").Value = 14000
$ws.Range("L3").Value = 14000
$ws.Range("M3").Value = 14000
$ws.Range("P3").Value = 560
$ws.Range("D4").Value = 44455
$ws.Range("K4").Value = 12000
$ws.Range("L4").Value = 12000
$ws.Range("M4").Value = 12000
$ws.Range("O4").Value = 'Región del Maule'
$ws.Range("P4").Value = 480
$ws.Range("D5").Value = 44383
$ws.Range("J5").Value = 120
$ws.Range("K5").Value = 12000
$ws.Range("L5").Value = 12000
$ws.Range("M5").Value = 12000
$ws.Range("P5").Value = 480
$ws.Range("D6").Value = 44398
$ws.Range("J6").Value = 200
$ws.Range("K6").Value = 15000
$ws.Range("L6").Value = 15000
$ws.Range("M6").Value = 15000
$ws.Range("O6").Value = 'Provincia del Elquí'
$ws.Range("P6").Value = 600
$ws.Range("D7").Value = 44441
$ws.Range("K7").Value = 10000
$ws.Range("L7").Value = 10000
$ws.Range("M7").Value = 10000
$ws.Range("O7").Value = 'Provincia del Elquí'
$ws.Range("P7").Value = 400
$ws.Range("D8").Value = 44462
$ws.Range("J8").Value = 300
$ws.Range("O8").Value = 'Región de O''Higgins'
$ws.Range("D9").Value = 44460
$ws.Range("J9").Value = 150
$ws.Range("D10").Value = 44167
$ws.Range("K10").Value = 8000
$ws.Range("L10").Value = 8000
$ws.Range("M10").Value = 8000
$ws.Range("O10").Value = 'Región del Maule'
$ws.Range("P10").Value = 320
$ws.Range("D11").Value = 44446
$ws.Range("J11").Value = 200
$ws.Range("K11").Value = 10000
$ws.Range("L11").Value = 10000
$ws.Range("M11").Value = 10000
$ws.Range("O11").Value = 'Provincia del Elquí'
$ws.Range("P11").Value = 400
$ws.Range("D12").Value = 44449
$ws.Range("J12").Value = 200
$ws.Range("K12").Value = 12000
$ws.Range("L12").Value = 12000
$ws.Range("M12").Value = 12000
$ws.Range("O12").Value = 'Provincia del Elquí'
$ws.Range("P12").Value = 480
$ws.Range("D13").Value = 44463
$ws.Range("J13").Value = 300
$ws.Range("K13").Value = 10000
$ws.Range("L13").Value = 10000
$ws.Range("M13").Value = 10000
$ws.Range("O13").Value = 'Región de O''Higgins'
$ws.Range("P13").Value = 400
$ws.Range("D14").Value = 44447
$ws.Range("J14").Value = 200
$ws.Range("D15").Value = 44169
$ws.Range("K15").Value = 9000
$ws.Range("L15").Value = 9000
$ws.Range("M15").Value = 9000
$ws.Range("O15").Value = 'Región del Maule'
$ws.Range("P15").Value = 360
$ws.Range("D16").Value = 44467
$ws.Range("J16").Value = 300
$ws.Range("O16").Value = 'Región de O''Higgins'
$ws.Range("D17").Value = 44473
$ws.Range("J17").Value = 500
$ws.Range("K17").Value = 9000
$ws.Range("L17").Value = 9000
$ws.Range("M17").Value = 9000
$ws.Range("O17").Value = 'Región de O''Higgins'
$ws.Range("P17").Value = 360
$ws.Range("D18").Value = 44469
$ws.Range("J18").Value = 500
$ws.Range("K18").Value = 9000
$ws.Range("L18").Value = 9000
$ws.Range("M18").Value = 9000
$ws.Range("O18").Value = 'Región de O''Higgins'
$ws.Range("P18").Value = 360
$ws.Range("D19").Value = 44168
$ws.Range("K19").Value = 9000
$ws.Range("L19").Value = 9000
$ws.Range("M19").Value = 9000
$ws.Range("P19").Value = 360
$ws.Range("D20").Value = 44161
$ws.Range("J20").Value = 200
$ws.Range("K20").Value = 9000
$ws.Range("L20").Value = 9000
$ws.Range("M20").Value = 9000
$ws.Range("P20").Value = 360
$ws.Range("D21").Value = 44445
$ws.Range("J21").Value = 200
$ws.Range("K21").Value = 10000
$ws.Range("L21").Value = 10000
$ws.Range("M21").Value = 10000
$ws.Range("P21").Value = 400
$ws.Range("D23").Value = 44475
$ws.Range("J23").Value = 400
$ws.Range("K23").Value = 8000
$ws.Range("L23").Value = 8000
$ws.Range("M23").Value = 8000
$ws.Range("P23").Value = 320
$ws.Range("D24").Value = 44474
$ws.Range("J24").Value = 500
$ws.Range("K24").Value = 8500
$ws.Range("L24").Value = 8500
$ws.Range("M24").Value = 8500
$ws.Range("O24").Value = 'Región de O''Higgins'
$ws.Range("P24").Value = 340
$ws.Range("D25").Value = 44448
$ws.Range("J25").Value = 150
$ws.Range("K25").Value = 13000
$ws.Range("L25").Value = 13000
$ws.Range("M25").Value = 13000
$ws.Range("O25").Value = 'Provincia del Elquí'
$ws.Range("P25").Value = 520
$ws.Range("D26").Value = 44452
$ws.Range("J26").Value = 200
$ws.Range("K26").Value = 12000
$ws.Range("L26").Value = 12000
$ws.Range("M26").Value = 12000
$ws.Range("O26").Value = 'Provincia del Elquí'
$ws.Range("P26").Value = 480
$ws.Range("D27").Value = 44453
$ws.Range("K27").Value = 12000
$ws.Range("L27").Value = 12000
$ws.Range("M27").Value = 12000
$ws.Range("O27").Value = 'Provincia del Elquí'
$ws.Range("P27").Value = 480
$ws.Range("D28").Value = 44159
$ws.Range("J28").Value = 300
$ws.Range("K28").Value = 7000
$ws.Range("L28").Value = 7000
$ws.Range("M28").Value = 7000
$ws.Range("P28").Value = 280
$ws.Range("D29").Value = 44166
$ws.Range("J29").Value = 200
$ws.Range("K29").Value = 8000
$ws.Range("L29").Value = 8000
$ws.Range("M29").Value = 8000
$ws.Range("O29").Value = 'Región del Maule'
$ws.Range("P29").Value = 320
$ws.Range("D30").Value = 44468
$ws.Range("J30").Value = 300
$ws.Range("K30").Value = 9000
$ws.Range("L30").Value = 9000
$ws.Range("M30").Value = 9000
$ws.Range("O30").Value = 'Región de O''Higgins'
$ws.Range("P30").Value = 360
$ws.Range("D31").Value = 44466
$ws.Range("J31").Value = 300
$ws.Range("O31").Value = 'Región de O''Higgins'

# --- Append new row 32 ---
$ws.Range("A32").Value = 5
$ws.Range("B32").Value = 'Macroferia Regional de Talca'
$ws.Range("C32").Value = 'Maule'
$ws.Range("D32").Value = 44461
$ws.Range("E32").Value = 7
$ws.Range("F32").Value = 100112026
$ws.Range("G32").Value = 'Haba'
$ws.Range("H32").Value = 'Sin especificar'
$ws.Range("I32").Value = 'Primera'
$ws.Range("J32").Value = 200
$ws.Range("K32").Value = 12000
$ws.Range("L32").Value = 12000
$ws.Range("M32").Value = 12000
$ws.Range("N32").Value = '$/saco 25 kilos'
$ws.Range("O32").Value = 'Provincia del Elquí'
$ws.Range("P32").Value = 480
$ws.Range("Q32").Value = 25
$ws.Range("R32").Value = 'Hortaliza'

# Match the date display format used by the rest of column D (style "s=2"
# in the OOXML, numFmt "YYYY-MM-DD HH:MM:SS").
$ws.Range("D32").NumberFormat = "YYYY-MM-DD HH:MM:SS"
